$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Coin / Link / Price / Volume(1h) values per row (row number => values).
# Only rows 45 and 46 change Coin/Link (two entries swapped position); all
# rows 2-51 get refreshed Price / Volume(1h) quotes.
$rows = @{
    2 = @($null, $null, "28.626.98", "  +1.53%  ")
    3 = @($null, $null, "1.867.60", "  +1.95%  ")
    4 = @($null, $null, "1.006", "  +0.07%  ")
    5 = @($null, $null, "326.90", "  -1.03%  ")
    6 = @($null, $null, "1.005", "  +0.19%  ")
    7 = @($null, $null, "0.4633", "  +0.98%  ")
    8 = @($null, $null, "0.3905", "  +1.54%  ")
    9 = @($null, $null, "0.07900", "  +0.81%  ")
    10 = @($null, $null, "0.9705", "  +0.73%  ")
    11 = @($null, $null, "22.28", "  +2.32%  ")
    12 = @($null, $null, "1.798.58", "  +0.99%  ")
    13 = @($null, $null, "5.724", "  +0.29%  ")
    14 = @($null, $null, "6.928", "  +0.52%  ")
    15 = @($null, $null, "0.06914", "  +0.69%  ")
    16 = @($null, $null, "88.55", "  +2.36%  ")
    17 = @($null, $null, "1.006", "  +0.04%  ")
    18 = @($null, $null, "0.00001003", "  +1.24%  ")
    19 = @($null, $null, "16.91", "  +0.61%  ")
    20 = @($null, $null, "1.004", "  +0.20%  ")
    21 = @($null, $null, "28.599.28", "  +1.30%  ")
    22 = @($null, $null, "5.326", "  +0.02%  ")
    23 = @($null, $null, "11.07", "  +0.44%  ")
    24 = @($null, $null, "2.128", "  -1.17%  ")
    25 = @($null, $null, "2.141.69", "  +6.40%  ")
    26 = @($null, $null, "155.11", "  +1.47%  ")
    27 = @($null, $null, "19.32", "  +0.74%  ")
    28 = @($null, $null, "5.734", "  -1.21%  ")
    29 = @($null, $null, "1.991", "  +1.59%  ")
    30 = @($null, $null, "119.25", "  +2.44%  ")
    31 = @($null, $null, "0.09336", "  +0.14%  ")
    32 = @($null, $null, "0.9378", "  +0.12%  ")
    33 = @($null, $null, "5.312", "  +0.72%  ")
    34 = @($null, $null, "1.336", "  +1.06%  ")
    35 = @($null, $null, "3.344", "  -3.08%  ")
    36 = @($null, $null, "0.05809", "  -3.92%  ")
    37 = @($null, $null, "0.02111", "  -2.02%  ")
    38 = @($null, $null, "1.157", "  +0.25%  ")
    39 = @($null, $null, "7.893", "  +4.93%  ")
    40 = @($null, $null, "0.5650", "  +0.94%  ")
    41 = @($null, $null, "9.934", "  -0.31%  ")
    42 = @($null, $null, "0.1773", "  -0.24%  ")
    43 = @($null, $null, "0.07271", "  +3.58%  ")
    44 = @($null, $null, "2.226", "  -1.75%  ")
    45 = @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "11.75", "  +1.07%  ")
    46 = @("Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.5321", "  +0.71%  ")
    47 = @($null, $null, "1.138", "  -8.41%  ")
    48 = @($null, $null, "1.847", "  +0.56%  ")
    49 = @($null, $null, "113.61", "  +0.74%  ")
    50 = @($null, $null, "2.349", "  +1.44%  ")
    51 = @($null, $null, "1.004", "  +0.23%  ")
}

foreach ($row in $rows.Keys) {
    $vals = $rows[$row]
    $coin = $vals[0]
    $link = $vals[1]
    $price = $vals[2]
    $volume = $vals[3]

    if ($null -ne $coin) { $ws.Range("B$row").Value = $coin }
    if ($null -ne $link) { $ws.Range("C$row").Value = $link }

    # Price strings look numeric ("1.006", "28.626.98", ...); force the cell
    # to Text before writing so Excel keeps the literal string instead of
    # coercing it to a number/date, then drop the format change again so the
    # cell keeps its original (default) style.
    $priceCell = $ws.Range("D$row")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.ClearFormats()

    $ws.Range("E$row").Value = $volume
}
